# Add 2024 combine stats for existing rookies + a new rookie row (Tahj Washington).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 ---
$ws.Range("K4").Value = 4.45
$ws.Range("L4").Value = 39
$ws.Range("N4").Formula = "=10*12+4"
$ws.Range("O4").Value = 6.88
$ws.Range("P4").Value = 4.03

# --- Row 5 ---
$ws.Range("K5").Value = 4.33
$ws.Range("L5").Value = 38.5
$ws.Range("M5").Value = 11
$ws.Range("N5").Value = 126

# --- Row 6 ---
$ws.Range("K6").Value = 4.41
$ws.Range("L6").Value = 39
$ws.Range("N6").Formula = "=10*12+4"
$ws.Range("O6").Value = 6.9
$ws.Range("P6").Value = 4.31

# --- Row 7 ---
$ws.Range("K7").Value = 4.39
$ws.Range("L7").Value = 36
$ws.Range("M7").Value = 13
$ws.Range("N7").Value = 124

# --- Row 8 ---
$ws.Range("K8").Value = 4.36
$ws.Range("L8").Value = 40.5
$ws.Range("N8").Formula = "=11*12+2"

# --- Row 9 ---
$ws.Range("K9").Value = 4.61
$ws.Range("L9").Value = 38
$ws.Range("N9").Formula = "=10*12+7"

# --- Row 10 ---
$ws.Range("K10").Value = 4.34
$ws.Range("L10").Value = 39.5
$ws.Range("N10").Formula = "=11*12+4"

# --- Row 11 ---
$ws.Range("K11").Value = 4.39
$ws.Range("L11").Value = 40.5
$ws.Range("N11").Value = 126

# --- Row 12 ---
$ws.Range("K12").Value = 4.45
$ws.Range("L12").Value = 38.5
$ws.Range("N12").Formula = "=11*12+1"

# --- Row 13 ---
$ws.Range("K13").Value = 4.21
$ws.Range("L13").Value = 41
$ws.Range("N13").Formula = "=10*12+11"

# --- Row 14 ---
$ws.Range("K14").Value = 4.52
$ws.Range("L14").Value = 37.5
$ws.Range("N14").Formula = "=10*12+9"

# --- Row 15 ---
$ws.Range("K15").Value = 4.5
$ws.Range("L15").Value = 36.5
$ws.Range("M15").Value = 13
$ws.Range("N15").Value = 119

# --- Row 16 ---
$ws.Range("K16").Value = 4.41
$ws.Range("L16").Value = 42
$ws.Range("M16").Value = 17
$ws.Range("N16").Formula = "=10*12+9"
$ws.Range("O16").Value = 6.64
$ws.Range("P16").Value = 4.05

# --- Row 17 ---
$ws.Range("K17").Value = 4.47
$ws.Range("L17").Value = 37
$ws.Range("N17").Formula = "=10*12+7"
$ws.Range("O17").Value = 6.94
$ws.Range("P17").Value = 4.18

# --- Row 18 ---
$ws.Range("K18").Value = 4.52
$ws.Range("L18").Value = 37
$ws.Range("N18").Formula = "=10*12+8"
$ws.Range("P18").Value = 4.11

# --- Row 19 ---
$ws.Range("K19").Value = 4.39
$ws.Range("M19").Value = 12

# --- Row 21 ---
$ws.Range("K21").Value = 4.46
$ws.Range("L21").Value = 34
$ws.Range("N21").Value = 120
$ws.Range("O21").Value = 7.16

# --- Row 23 ---
$ws.Range("M23").Value = 21

# --- Row 26 ---
$ws.Range("K26").Value = 4.47
$ws.Range("L26").Value = 42.5
$ws.Range("M26").Value = 19
$ws.Range("N26").Formula = "=10*12+6"

# --- Row 27 ---
$ws.Range("K27").Value = 4.54
$ws.Range("L27").Value = 37
$ws.Range("N27").Value = 121

# --- Row 28 ---
$ws.Range("K28").Value = 4.57
$ws.Range("L28").Value = 42.5
$ws.Range("N28").Value = 128

# --- Row 29 ---
$ws.Range("K29").Value = 4.38
$ws.Range("L29").Value = 36
$ws.Range("N29").Formula = "=9*12+11"
$ws.Range("O29").Value = 7.02
$ws.Range("P29").Value = 4.32

# --- Row 30 ---
$ws.Range("K30").Value = 4.44
$ws.Range("L30").Value = 37.5
$ws.Range("N30").Formula = "=10*12+7"

# --- Row 31 ---
$ws.Range("K31").Value = 4.46
$ws.Range("L31").Value = 36
$ws.Range("N31").Value = 121
$ws.Range("O31").Value = 6.7
$ws.Range("P31").Value = 4.02

# --- New row 34: Tahj Washington (2024) ---
$ws.Range("A34").Value = "Tahj Washington"
$ws.Range("B34").Value = 2024

# Copy the K:P number formatting/style from an existing data row so the new
# blank stat cells pick up style s="2" like the rest of the sheet.
$ws.Range("K4:P4").Copy()
$ws.Range("K34:P34").PasteSpecial(-4122)

$ws.Range("L34").Value = 35
$ws.Range("N34").Formula = "=10*12+2"

$ws.Range("U34").Value = 0
$ws.Range("V34").Value = 0
$ws.Range("W34").Value = 0
$ws.Range("X34").Value = 0
$ws.Range("Y34").Formula = "=YEAR(TODAY())-B34"
